$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# --- Overview sheet: Status column for both rows becomes "Handed back..." ---
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("B2").Value = $statusHandedBack
$ov.Range("C2").Value = $statusHandedBack
$ov.Range("B3").Value = $statusHandedBack
$ov.Range("C3").Value = $statusHandedBack

# --- zh-cn sheet ---
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("B2").Value = $statusHandedBack
$zh.Range("E2").Value = "25f2e471-237f-4a05-95e5-5daa681a7bc7.md"
$zh.Range("F2").Value = "25f2e471-237f-4a05-95e5-5daa681a7bc7.9081ab96605150e58cb69fe39e825b77e885e6cf.zh-cn.xlf"
$zh.Range("G2").Value = "2016-02-22 18:08:22"

$zh.Range("B3").Value = $statusHandedBack
$zh.Range("E3").Value = "df44ebf6-6bde-4013-801a-af2484337ce6.md"
$zh.Range("F3").Value = "df44ebf6-6bde-4013-801a-af2484337ce6.f400ca1a38bfa6ac5697dd326043e7ba7ab0b971.zh-cn.xlf"
$zh.Range("G3").Value = "2016-02-22 18:08:22"

# --- de-de sheet ---
$de = $wb.Worksheets.Item("de-de")
$de.Range("B2").Value = $statusHandedBack
$de.Range("E2").Value = "25f2e471-237f-4a05-95e5-5daa681a7bc7.md"
$de.Range("F2").Value = "25f2e471-237f-4a05-95e5-5daa681a7bc7.9081ab96605150e58cb69fe39e825b77e885e6cf.de-de.xlf"
$de.Range("G2").Value = "2016-02-22 18:08:42"

$de.Range("B3").Value = $statusHandedBack
$de.Range("E3").Value = "df44ebf6-6bde-4013-801a-af2484337ce6.md"
$de.Range("F3").Value = "df44ebf6-6bde-4013-801a-af2484337ce6.f400ca1a38bfa6ac5697dd326043e7ba7ab0b971.de-de.xlf"
$de.Range("G3").Value = "2016-02-22 18:08:42"
